# Fruta / hortaliza, semanal
# Insert two new weekly rows (Primera / Segunda) at the top of the
# "Betarraga" data block, pushing the existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 135:170 down to 137:172, leaving two blank rows behind.
$ws.Rows("135:136").Insert()

# --- New row 135 (Primera) ---
$ws.Range("A135").Value = 1
$ws.Range("B135").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C135").Value = "Arica y Parinacota"
$ws.Range("D135").Value = 44476
$ws.Range("E135").Value = 15
$ws.Range("F135").Value = 100114014
$ws.Range("G135").Value = "Betarraga"
$ws.Range("H135").Value = "Sin especificar"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 400
$ws.Range("L135").Value = 450
$ws.Range("M135").Value = 425
$ws.Range("N135").Value = "`$/paquete 4 unidades"
$ws.Range("O135").Value = "Región de Arica y Parinacota"
$ws.Range("P135").Value = 106
$ws.Range("Q135").Value = 4
$ws.Range("R135").Value = "Hortaliza"

# --- New row 136 (Segunda) ---
$ws.Range("A136").Value = 1
$ws.Range("B136").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C136").Value = "Arica y Parinacota"
$ws.Range("D136").Value = 44476
$ws.Range("E136").Value = 15
$ws.Range("F136").Value = 100114014
$ws.Range("G136").Value = "Betarraga"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Segunda"
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 400
$ws.Range("L136").Value = 450
$ws.Range("M136").Value = 425
$ws.Range("N136").Value = "`$/paquete 5 unidades"
$ws.Range("O136").Value = "Región de Arica y Parinacota"
$ws.Range("P136").Value = 85
$ws.Range("Q136").Value = 5
$ws.Range("R136").Value = "Hortaliza"
